$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(5)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Text edits (applied right-to-left so earlier offsets stay valid) ---

# 3) ", no caso as diferentes implementações que utilizam objetos do tipo Actions como botões e menus."
#    -> ", no caso podem ser criados menus e botões que podem realizar diversas ações.  "
$tr.Characters(333, 96).Text = ", no caso podem ser criados menus e botões que podem realizar diversas ações.  "

# 2) "), aqui a classe criadora" -> "
#        ) visto que esta é comum a todos os objetos criados e os outros parâmetros
#        passados são classes java como o " + "command" (own run) + ", aqui a classe criadora"
$tr.Characters(201, 25).Text = ") visto que esta é comum a todos os objetos criados e os outros parâmetros passados são classes java como o command, aqui a classe criadora"

# Split "command" into its own run (re-apply the same size so formatting is untouched
# but a run boundary is forced around just that word).
$full = $tr.Text
$cmdStart = $full.IndexOf("como o command") + 1 + ("como o ").Length
$cmdRange = $tr.Characters($cmdStart, 7)
$cmdRange.Font.Size = 12

# 1) "é a classe criadora, o produto é a interface " -> "é a classe criadora, o produto pode ser considerado a interface "
$tr.Characters(126, 45).Text = "é a classe criadora, o produto pode ser considerado a interface "

# --- Resize/move the textbox to fit the new (longer) text ---
$sh.Top = [double]215.812
$sh.Height = [double]96.513386
